$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the single pre-existing hyperlink (on the Password cell, D3) so it can
# be re-added cleanly alongside the new one below.
$ws.Hyperlinks.Delete()

# Password value changes (R@ckware -> R@ckware4IT)
$ws.Range("D3").Value = "R@ckware4IT"

# New validation: Username (C3) becomes a mailto hyperlink, same pattern
# as the existing Password hyperlink.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:administrator@vsphere.local")
$ws.Range("C3").Style = "Hyperlink"

# Re-add the Password hyperlink pointing at the new value.
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:R@ckware4IT")
$ws.Range("D3").Style = "Hyperlink"

# Port is now a known value (443) instead of "NA"
$ws.Range("E3").Value = 443

# Move the active selection (cosmetic, matches recorded cursor position)
$ws.Range("G8").Select() | Out-Null
